$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# MyersBriggs S/T/P counts (B3:B5) were recomputed
$ws.Range("B3").Value = 10
$ws.Range("B4").Value = 9
$ws.Range("B5").Value = 17

# Row 2 (header detail row) got shorter after the note row beneath it was removed
$ws.Rows(2).RowHeight = 187.2

# Remove the note row ("This Strengths order matches the order on the team
# CliftonStrengths chart") that used to sit at B7, between the summary rows
# and the Strengths list.
$ws.Range("B7").ClearContents()

# Re-sort the Strengths list (B8:B41) alphabetically instead of grouped by
# CliftonStrengths domain, and persist the sort state like Excel's Sort
# dialog does.
$sortRange = $ws.Range("B8:B41")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Update the view: scrolled down a bit and the active selection moved to B7.
$ws.Range("B7").Select()
